$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 6, shifting rows 6-8 down to 7-9
$ws.Rows.Item(6).Insert()

# Fill in the new row 6 with the new data
$ws.Cells.Item(6, 1).Value = "Torneo FEG"
$ws.Cells.Item(6, 2).Value = "Principiantes"
$ws.Cells.Item(6, 3).Value = "general"
$ws.Cells.Item(6, 4).Value = 5
$ws.Cells.Item(6, 5).Value = "Bogado, Ogán"
$ws.Cells.Item(6, 6).Value = 33
